$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 41; Id = 40; Categoria = "Camicado"; Titulo = "Colocação de papel de parede"; Preco = 200.0012; Imagem = "https://img.camicado.com.br/item/sku580002/small/1.jpg"; Disponivel = 2 },
    @{ Row = 42; Id = 41; Categoria = "Camicado"; Titulo = "Colocação de piso"; Preco = 300.0013; Imagem = "https://img.camicado.com.br/item/sku580003/small/1.jpg"; Disponivel = 3 },
    @{ Row = 43; Id = 42; Categoria = "Camicado"; Titulo = "Material para reforma de ambiente"; Preco = 400.0014; Imagem = "https://img.camicado.com.br/item/sku580004/small/1.jpg"; Disponivel = 4 },
    @{ Row = 44; Id = 43; Categoria = "Camicado"; Titulo = "Instalação de portas"; Preco = 300.0013; Imagem = "https://img.camicado.com.br/item/sku580012/small/1.jpg"; Disponivel = 5 },
    @{ Row = 45; Id = 44; Categoria = "Camicado"; Titulo = "Conjunto malas para viagem"; Preco = 450.00145; Imagem = "https://img.camicado.com.br/item/sku580016/small/1.jpg"; Disponivel = 6 },
    @{ Row = 46; Id = 45; Categoria = "Camicado"; Titulo = "Prestador de serviços residenciais"; Preco = 200.0012; Imagem = "https://img.camicado.com.br/item/sku580018/small/1.jpg"; Disponivel = 7 },
    @{ Row = 47; Id = 46; Categoria = "Camicado"; Titulo = "Pintura ambiente"; Preco = 300.0013; Imagem = "https://img.camicado.com.br/item/sku580019/small/1.jpg"; Disponivel = 8 },
    @{ Row = 48; Id = 47; Categoria = "Camicado"; Titulo = "Aula de Valsa para os Noivos"; Preco = 180.00118; Imagem = "https://img.camicado.com.br/item/sku580029/small/1.jpg"; Disponivel = 9 },
    @{ Row = 49; Id = 48; Categoria = "Camicado"; Titulo = "Dia da noiva"; Preco = 800.0018; Imagem = "https://img.camicado.com.br/item/sku580041/small/1.jpg"; Disponivel = 10 },
    @{ Row = 50; Id = 49; Categoria = "Camicado"; Titulo = "Massagem com pedras quentes"; Preco = 120.00112; Imagem = "https://img.camicado.com.br/item/sku580050/small/1.jpg"; Disponivel = 11 },
    @{ Row = 51; Id = 50; Categoria = "Camicado"; Titulo = "Massagem relaxante"; Preco = 90.0019; Imagem = "https://img.camicado.com.br/item/sku580052/small/1.jpg"; Disponivel = 12 },
    @{ Row = 52; Id = 51; Categoria = "Camicado"; Titulo = "Diária em Hotel"; Preco = 180.00118; Imagem = "https://img.camicado.com.br/item/sku580059/small/1.jpg"; Disponivel = 13 },
    @{ Row = 53; Id = 52; Categoria = "Camicado"; Titulo = "Diária em pousada no campo"; Preco = 150.00115; Imagem = "https://img.camicado.com.br/item/sku580062/small/1.jpg"; Disponivel = 14 },
    @{ Row = 54; Id = 53; Categoria = "Camicado"; Titulo = "Transporte para Mudança"; Preco = 300.0013; Imagem = "https://img.camicado.com.br/item/sku580079/small/1.jpg"; Disponivel = 15 },
    @{ Row = 55; Id = 54; Categoria = "Camicado"; Titulo = "Aula de Culinária"; Preco = 100.0011; Imagem = "https://img.camicado.com.br/item/sku580092/small/1.jpg"; Disponivel = 16 },
    @{ Row = 56; Id = 55; Categoria = "Camicado"; Titulo = "Enxoval de lingerie"; Preco = 400.0014; Imagem = "https://img.camicado.com.br/item/sku580105/small/1.jpg"; Disponivel = 17 },
    @{ Row = 57; Id = 56; Categoria = "Camicado"; Titulo = "Almoço no centro da cidade"; Preco = 140.00114; Imagem = "https://img.camicado.com.br/item/sku580138/small/1.jpg"; Disponivel = 18 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Id
    $ws.Cells.Item($r.Row, 2).Value = "disponivel"
    $ws.Cells.Item($r.Row, 3).Value = $r.Categoria
    $ws.Cells.Item($r.Row, 4).Value = $r.Titulo
    $ws.Cells.Item($r.Row, 6).Value = $r.Preco
    $ws.Cells.Item($r.Row, 8).Value = $r.Imagem
    $ws.Cells.Item($r.Row, 9).Value = $r.Disponivel
}

Write-Output "done"